$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new column before D, shifting existing D:K data to E:L
$ws.Columns("D:D").Insert()

# Propagate number formatting from (now-shifted) column E into the new column D
$ws.Range("E7:E102").Copy()
$ws.Range("D7:D102").PasteSpecial(-4122)

# Populate the new column D plus the handful of cells whose values changed
# beyond a simple shift (re-stated figures for FY2018 reporting)
$ws.Range("D7").Value = 43465
$ws.Range("D8").Value = 4979000
$ws.Range("E8").Value = 4638000
$ws.Range("F8").Value = 4512000
$ws.Range("D9").Value = "NA"
$ws.Range("D10").Value = "NA"
$ws.Range("D12").Value = "NA"
$ws.Range("D13").Value = 0
$ws.Range("D14").Value = 30000
$ws.Range("D15").Value = 586000
$ws.Range("D17").Value = 2392000
$ws.Range("E17").Value = 2259000
$ws.Range("F17").Value = 2340000
$ws.Range("D18").Value = 2587000
$ws.Range("F18").Value = 2172000
$ws.Range("D20").Value = 181000
$ws.Range("E20").Value = 334000
$ws.Range("F20").Value = 49000
$ws.Range("D21").Value = 3354000
$ws.Range("E21").Value = 3248000
$ws.Range("F21").Value = 2831000
$ws.Range("D22").Value = 244000
$ws.Range("D23").Value = 2524000
$ws.Range("E23").Value = 2526000
$ws.Range("F23").Value = 2043000
$ws.Range("D24").Value = 511000
$ws.Range("E24").Value = 736000
$ws.Range("F24").Value = 586000
$ws.Range("D25").Value = 0
$ws.Range("D26").Value = 2013000
$ws.Range("E26").Value = 1790000
$ws.Range("F26").Value = 1457000
$ws.Range("D27").Value = 1981000
$ws.Range("E27").Value = 1762000
$ws.Range("F27").Value = 1430000
$ws.Range("D28").Value = 0
$ws.Range("D29").Value = 11000
$ws.Range("D30").Value = 0
$ws.Range("D31").Value = 0
$ws.Range("D32").Value = -181000
$ws.Range("E32").Value = -334000
$ws.Range("F32").Value = -49000
$ws.Range("D33").Value = 1992000
$ws.Range("E33").Value = 2526000
$ws.Range("F33").Value = 1430000
$ws.Range("D34").Value = 0
$ws.Range("D35").Value = 1992000
$ws.Range("E35").Value = 2526000
$ws.Range("F35").Value = 1430000
$ws.Range("D38").Value = 43465
$ws.Range("D41").Value = 724000
$ws.Range("D42").Value = "NA"
$ws.Range("D43").Value = 953000
$ws.Range("D44").Value = 0
$ws.Range("D45").Value = 65015000
$ws.Range("D46").Value = 66692000
$ws.Range("D47").Value = "NA"
$ws.Range("D48").Value = 1241000
$ws.Range("D49").Value = 23547000
$ws.Range("D50").Value = 0
$ws.Range("D51").Value = 0
$ws.Range("D52").Value = 1311000
$ws.Range("D53").Value = 0
$ws.Range("D54").Value = 92791000
$ws.Range("D57").Value = 521000
$ws.Range("D58").Value = 951000
$ws.Range("D59").Value = 64636000
$ws.Range("D60").Value = 66108000
$ws.Range("D61").Value = 6490000
$ws.Range("D62").Value = 2891000
$ws.Range("D63").Value = 0
$ws.Range("D64").Value = 0
$ws.Range("D65").Value = 0
$ws.Range("D66").Value = 75590000
$ws.Range("D68").Value = 0
$ws.Range("D69").Value = 0
$ws.Range("D70").Value = 0
$ws.Range("D71").Value = 0
$ws.Range("D72").Value = 8317000
$ws.Range("D73").Value = 0
$ws.Range("D74").Value = 0
$ws.Range("D75").Value = 0
$ws.Range("D76").Value = 17201000
$ws.Range("D77").Value = 0
$ws.Range("D80").Value = 43465
$ws.Range("D81").Value = 1992000
$ws.Range("E81").Value = 2526000
$ws.Range("F81").Value = 1430000
$ws.Range("D83").Value = 586000
$ws.Range("D84").Value = 0
$ws.Range("D85").Value = 0
$ws.Range("D86").Value = 0
$ws.Range("D87").Value = 0
$ws.Range("D88").Value = 0
$ws.Range("D89").Value = 2533000
$ws.Range("D91").Value = -134000
$ws.Range("E91").Value = -220000
$ws.Range("F91").Value = -250000
$ws.Range("G91").Value = -190000
$ws.Range("H91").Value = -172000
$ws.Range("I91").Value = -134000
$ws.Range("J91").Value = -32400
$ws.Range("D92").Value = 0
$ws.Range("D93").Value = 0
$ws.Range("D94").Value = -1755000
$ws.Range("D96").Value = -555000
$ws.Range("D97").Value = 0
$ws.Range("D98").Value = 0
$ws.Range("D99").Value = 0
$ws.Range("D100").Value = -463000
$ws.Range("D101").Value = -11000
$ws.Range("D102").Value = 304000

$excel.CutCopyMode = 0
